# ---------------------------------------------------------------------------
# res_stock_pnw_existing.xlsx edit:
#  - Bump OpenStudio Server Version 1.18.0-rc0 -> 1.19.0-rc0 (Setup!B5)
#  - Worker Nodes 3 -> 0 (Setup!B9, formula in E9 recalculates automatically)
#  - Number of Samples 10000 -> 100 (Setup!B24 / D24)
#  - Add 5 new "airflow" measures to the Variables sheet (Ducts, Infiltration,
#    Natural Ventilation, Mechanical Ventilation, Airflow -- Airflow left
#    disabled because of an OS error, per the commit message)
#  - Add matching 5 new rows to the Outputs sheet for the same measures
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Setup sheet
# ---------------------------------------------------------------------------
$setup = $wb.Worksheets.Item("Setup")
$setup.Range("B5").Value = "1.19.0-rc0"
$setup.Range("B9").Value = 0
# B24 is a formula (=IF(D24<>"",D24,...)); only the literal D24 input changes,
# the formula result recalculates to the same 100 automatically.
$setup.Range("D24").Value = 100

# ---------------------------------------------------------------------------
# 2. Variables sheet - insert 5 new measure blocks (3 rows each) right before
#    the existing reporting-measures block that starts at row 192.
# ---------------------------------------------------------------------------
$vars = $wb.Worksheets.Item("Variables")

# Make room: 5 measures x 3 rows = 15 rows, inserted above row 192.
$vars.Range("192:206").Insert()

# Use the existing "Cooling/Heating Setpoint" style measure block (rows
# 141:143) as a formatting template -- copy it five times into the newly
# opened rows, then overwrite the cell values per measure below.
$template = $vars.Range("A141:X143")
$template.Copy($vars.Range("A192:X194"))
$template.Copy($vars.Range("A195:X197"))
$template.Copy($vars.Range("A198:X200"))
$template.Copy($vars.Range("A201:X203"))
$template.Copy($vars.Range("A204:X206"))

# NOTE: this PowerShell host does not bind named (-Param value) arguments,
# so helper functions below use positional parameters only.
function Set-AirflowMeasure {
    param($HeaderRow, $Enabled, $SetName, $TsvName, $SampleValueName)

    $r1 = $HeaderRow
    $r2 = $HeaderRow + 1
    $r3 = $HeaderRow + 2

    $vars.Range("A$r1").Value = $Enabled
    $vars.Range("B$r1").Value = $SetName

    $vars.Range("I$r2").Value = $TsvName

    $vars.Range("D$r3").Value = $SampleValueName
}

Set-AirflowMeasure 192 $true  "Set Ducts"                  "Ducts.tsv"                  "Ducts Sample Value"
Set-AirflowMeasure 195 $true  "Set Infiltration"           "Infiltration.tsv"           "Infiltration Sample Value"
Set-AirflowMeasure 198 $true  "Set Natural Ventilation"    "Natural Ventilation.tsv"    "Natural Ventilation Sample Value"
Set-AirflowMeasure 201 $true  "Set Mechanical Ventilation" "Mechanical Ventilation.tsv" "Mechanical Ventilation Sample Value"
Set-AirflowMeasure 204 $false "Set Airflow"                "Airflow.tsv"                "Airflow Sample Value"

# ---------------------------------------------------------------------------
# 3. Outputs sheet - insert 5 new rows right before the existing row 65
#    ("Heating Setpoint" reporting output), one per new measure.
# ---------------------------------------------------------------------------
$outs = $wb.Worksheets.Item("Outputs")

$outs.Range("65:69").Insert()

$outTemplate = $outs.Range("A33:I33")
$outTemplate.Copy($outs.Range("A65:I65"))
$outTemplate.Copy($outs.Range("A66:I66"))
$outTemplate.Copy($outs.Range("A67:I67"))
$outTemplate.Copy($outs.Range("A68:I68"))
$outTemplate.Copy($outs.Range("A69:I69"))

function Set-AirflowOutput {
    param($Row, $Name, $ReportingName)
    $outs.Range("A$Row").Value = $Name
    $outs.Range("D$Row").Value = $ReportingName
}

Set-AirflowOutput 65 "Ducts"                  "res_stock_reporting.Ducts"
Set-AirflowOutput 66 "Infiltration"           "res_stock_reporting.Infiltration"
Set-AirflowOutput 67 "Natural Ventilation"    "res_stock_reporting.Natural Ventilation"
Set-AirflowOutput 68 "Mechanical Ventilation" "res_stock_reporting.Mechanical Ventilation"
Set-AirflowOutput 69 "Airflow"                "res_stock_reporting.Airflow"
